# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the data source.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 61
$wsExhibit.Range("F7").Value = 9648
$wsExhibit.Range("F9").Value = 331
$wsExhibit.Range("F10").Value = 1221
$wsExhibit.Range("F11").Value = 2765
$wsExhibit.Range("F15").Value = 24
$wsExhibit.Range("F17").Value = 478
$wsExhibit.Range("F18").Value = 100
$wsExhibit.Range("F19").Value = 260
$wsExhibit.Range("F20").Value = 1362

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 61
$wsAll.Range("F8").Value = 9648
$wsAll.Range("F10").Value = 331
$wsAll.Range("F11").Value = 1221
$wsAll.Range("F12").Value = 2765
$wsAll.Range("F16").Value = 24
$wsAll.Range("F18").Value = 478
$wsAll.Range("F19").Value = 100
$wsAll.Range("F20").Value = 260
$wsAll.Range("F21").Value = 1362
